# Updated test data for FC test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# K7 stays "MinBatteryLoadingDetail" - already set, but set explicitly to keep parity
$ws.Range("K7").Value = "MinBatteryLoadingDetail"

# K8 text changes from "Minimum battery (Ah)" to "Minimum Battery size(Ah)"
$ws.Range("K8").Value = "Minimum Battery size(Ah)"

# Reflect the final active selection on the sheet (K8)
$ws.Range("K8").Select()
